$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 4: new contact "alekhya" replaces "komal mam" ---
$ws.Range("B4").Value = "alekhya"
$ws.Range("C4").Value = "alekhyakanjarla@gmail.com"

# --- Insert new row 5: new contact "nishita mam" ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "nishita mam"
$ws.Range("C5").Value = "gojo.testing123@gmail.com"

# --- Rebuild hyperlinks (this runtime manages Hyperlinks per-sheet, not per-cell) ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:gojo.testing123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:guru.sai.shreesh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:gojo.testing123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:alekhyakanjarla@gmail.com")

# The new C4 hyperlink carries a cached "display" label left over from an earlier edit
$hl = $ws.Hyperlinks.Item($ws.Hyperlinks.Count)
$hl.TextToDisplay = "gojo.testing123@gmail.com"
$ws.Range("C4").Value = "alekhyakanjarla@gmail.com"

# Reapply the Hyperlink style so every linked cell in the column matches
$ws.Range("C2:C5").Style = "Hyperlink"

# Final selection lands on the newly entered cell
$null = $ws.Range("C5").Select()
